$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.202.95"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.474.17"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.27"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.70"
$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "2.475.48"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.68"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "2.914.76"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "58.304.96"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.09"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "2.472.11"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.44"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.16"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.50"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("E22").Value = "  +6.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.95"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "2.549.80"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  +7.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.75"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Value = "0.0₃0744"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.29"
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.14"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +4.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.828"
$ws.Range("E41").Value = "  +5.22%  "
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "36.58"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "137.84"
$ws.Range("E44").Value = "  +13.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.44"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.98"
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "263.76"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.578"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0218"
$ws.Range("E51").Value = "  +3.53%  "
